# Rename sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "species"
$ws2.Name = "turtles"

# Add data rows to sheet1 ("species")
$ws1.Range("A2").Value = "Dog"
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 36
$ws1.Range("D2").Value = 120

$ws1.Range("A3").Value = "Cat"
$ws1.Range("B3").Value = 2
$ws1.Range("C3").Value = 36
$ws1.Range("D3").Value = 120

# Add data row to sheet2 ("turtles")
$ws2.Range("A2").Value = "Big turtle"
$ws2.Range("B2").Value = 4

# Make the second sheet (turtles) the active one
$ws2.Activate()

# Set selections to match target state
$ws1.Range("D8").Select()
$ws2.Range("F10").Select()
$ws2.Activate()
